$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "weight distance" value (C25) - this also causes C26/C27 to recalc
$ws.Range("C25").Value = 0.14399999999999999

# Normalize style of B22/C22 and C25 to match the rest of the sheet (style used by C26/B26 etc.)
$ws.Range("B22").Style = $ws.Range("B21").Style
$ws.Range("C22").Style = $ws.Range("C26").Style
$ws.Range("C25").Style = $ws.Range("C26").Style
